$wb = $excel.ActiveWorkbook

# The "Generate Report for Handback" run discovered that the handback file for
# b7cd9a63-107a-43a2-8dd7-f267f7704a34 was stale, so row 7 on both the zh-cn
# and de-de sheets now records the (out of date) target file that was handed
# back, together with a new handback datetime and an error message explaining
# the version mismatch.

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---- zh-cn sheet (row 7) ----------------------------------------------
$zhcn.Range("I7").Value = "b7cd9a63-107a-43a2-8dd7-f267f7704a34.md"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267acef0fa4050bb66372fdc75e5f765b2978bd9/e2e/b7cd9a63-107a-43a2-8dd7-f267f7704a34.md",
    "",
    "",
    "b7cd9a63-107a-43a2-8dd7-f267f7704a34.md"
) | Out-Null

$zhcn.Range("J7").Value = "b7cd9a63-107a-43a2-8dd7-f267f7704a34.f500811b1c85ae8831eb860bbe5a5d1d86ff19ef.zh-cn.xlf"
$zhcn.Range("K7").Value = "2016-08-24 18:52:18"
$zhcn.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d1d7defebcf52d4e9f5c82f99afce8f4d6ba1c5/e2e/b7cd9a63-107a-43a2-8dd7-f267f7704a34.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267acef0fa4050bb66372fdc75e5f765b2978bd9/e2e/b7cd9a63-107a-43a2-8dd7-f267f7704a34.md."

# Widen the "Error Detail" column so the new message is readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet (row 7) -----------------------------------------------
$dede.Range("I7").Value = "b7cd9a63-107a-43a2-8dd7-f267f7704a34.md"
$dede.Hyperlinks.Add(
    $dede.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267acef0fa4050bb66372fdc75e5f765b2978bd9/e2e/b7cd9a63-107a-43a2-8dd7-f267f7704a34.md",
    "",
    "",
    "b7cd9a63-107a-43a2-8dd7-f267f7704a34.md"
) | Out-Null

$dede.Range("J7").Value = "b7cd9a63-107a-43a2-8dd7-f267f7704a34.f500811b1c85ae8831eb860bbe5a5d1d86ff19ef.de-de.xlf"
$dede.Range("K7").Value = "2016-08-24 18:52:26"
$dede.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d1d7defebcf52d4e9f5c82f99afce8f4d6ba1c5/e2e/b7cd9a63-107a-43a2-8dd7-f267f7704a34.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267acef0fa4050bb66372fdc75e5f765b2978bd9/e2e/b7cd9a63-107a-43a2-8dd7-f267f7704a34.md."

# Widen the "Error Detail" column so the new message is readable.
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

Write-Host "Report generated for handback b7cd9a63-107a-43a2-8dd7-f267f7704a34"
